$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number need to be forced back to
# text (the source workbook stores every "Price" figure as an inline string,
# e.g. "1.00" / "0.420" -- letting Excel auto-coerce them to numbers would
# drop the trailing zeros / change the stored value).
function Set-TextValue($cell, $val) {
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range("D2").Value = '59.769.43'
$ws.Range("E2").Value = '  +2.60%  '
$ws.Range("D3").Value = '2.422.62'
$ws.Range("E3").Value = '  +3.06%  '
$ws.Range("E4").Value = '  -0.06%  '
Set-TextValue "D5" '554.13'
$ws.Range("E5").Value = '  +2.41%  '
Set-TextValue "D6" '137.78'
$ws.Range("E6").Value = '  +1.96%  '
Set-TextValue "D7" '1.00'
$ws.Range("E7").Value = '  -0.04%  '
Set-TextValue "D8" '0.569'
$ws.Range("E8").Value = '  +1.46%  '
$ws.Range("E9").Value = '  +5.05%  '
$ws.Range("E10").Value = '  +3.36%  '
$ws.Range("E11").Value = '  +1.57%  '
$ws.Range("E12").Value = '  -2.04%  '
Set-TextValue "D13" '24.67'
$ws.Range("E13").Value = '  +3.34%  '
$ws.Range("E14").Value = '  +2.93%  '
$ws.Range("D15").Value = '59.655.59'
$ws.Range("E15").Value = '  +2.51%  '
$ws.Range("E16").Value = '  +4.59%  '
$ws.Range("D17").Value = '2.398.57'
$ws.Range("E17").Value = '  +1.91%  '
Set-TextValue "D18" '11.34'
$ws.Range("E18").Value = '  +5.89%  '
$ws.Range("E19").Value = '  +4.58%  '
Set-TextValue "D20" '334.97'
$ws.Range("E20").Value = '  +0.80%  '
Set-TextValue "D21" '6.94'
$ws.Range("E21").Value = '  +3.46%  '
Set-TextValue "D22" '1.00'
$ws.Range("E22").Value = '  +0.07%  '
Set-TextValue "D23" '64.63'
$ws.Range("E23").Value = '  +2.93%  '
$ws.Range("E24").Value = '  +0.97%  '
Set-TextValue "D25" '8.66'
$ws.Range("E25").Value = '  +2.26%  '
$ws.Range("E27").Value = '  -1.21%  '
$ws.Range("D28").Value = '0.0₃0790'
$ws.Range("E28").Value = '  +7.41%  '
$ws.Range("E29").Value = '  +2.95%  '
Set-TextValue "D30" '170.41'
$ws.Range("E30").Value = '  +0.36%  '
$ws.Range("E31").Value = '  +3.07%  '
Set-TextValue "D32" '18.73'
$ws.Range("E32").Value = '  +1.94%  '
$ws.Range("E33").Value = '  +1.15%  '
$ws.Range("E34").Value = '  -0.01%  '
Set-TextValue "D35" '1.31'
$ws.Range("E35").Value = '  +5.09%  '
Set-TextValue "D36" '4.27'
$ws.Range("E36").Value = '  +0.71%  '
$ws.Range("E37").Value = '  -0.14%  '
$ws.Range("E38").Value = '  -0.77%  '
Set-TextValue "D39" '40.11'
$ws.Range("E39").Value = '  +2.38%  '
Set-TextValue "D40" '0.420'
$ws.Range("E40").Value = '  +11.09%  '
Set-TextValue "D41" '313.56'
$ws.Range("E42").Value = '  +3.01%  '
Set-TextValue "D43" '142.85'
$ws.Range("E43").Value = '  -0.10%  '
$ws.Range("E44").Value = '  +2.66%  '
$ws.Range("B46").Value = 'Polygon'
$ws.Range("C46").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue "D46" '0.412'
$ws.Range("E46").Value = '  +6.90%  '
$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue "D47" '19.24'
$ws.Range("E47").Value = '  +0.38%  '
$ws.Range("B48").Value = 'Mantle'
$ws.Range("C48").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue "D48" '0.573'
$ws.Range("E48").Value = '  +1.60%  '
$ws.Range("E49").Value = '  +3.10%  '
Set-TextValue "D50" '11.04'
$ws.Range("E50").Value = '  -0.32%  '
$ws.Range("E51").Value = '  +5.13%  '
